$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 231.9
$ws.Range("B3").Value = 146
$ws.Range("C3").Value = 354.7
$ws.Range("C4").Value = 296.1
